# Fixes for VRelay and Upgrade Verification SCript
#
# Updates the Katalon bootstrap results workbook: flips the latest
# execution results from "Fail" to "Pass" and stamps the corresponding
# execution-date cells with the new run timestamps, across all three
# worksheets (CreateModifyDeleteProfile, AddModifyDeleteCC,
# AddModifyDeleteACH).

$wb = $excel.ActiveWorkbook

# --- CreateModifyDeleteProfile ---
$ws1 = $wb.Worksheets.Item("CreateModifyDeleteProfile")
$ws1.Range("A2").Value = "Pass"
$ws1.Range("B2").Value = "Fri Aug 22 23:06:45 IST 2025"

# --- AddModifyDeleteCC ---
$ws2 = $wb.Worksheets.Item("AddModifyDeleteCC")
$ws2.Range("A2").Value = "Pass"
$ws2.Range("B2").Value = "Fri Aug 22 23:05:31 IST 2025"

# --- AddModifyDeleteACH ---
$ws3 = $wb.Worksheets.Item("AddModifyDeleteACH")
$ws3.Range("A2").Value = "Pass"
$ws3.Range("B2").Value = "Fri Aug 22 23:01:18 IST 2025"

$ws3.Range("A3").Value = "Pass"
$ws3.Range("B3").Value = "Fri Aug 22 23:02:39 IST 2025"

$ws3.Range("A4").Value = "Pass"
$ws3.Range("B4").Value = "Fri Aug 22 23:04:02 IST 2025"
